$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 72; this shifts the previous
# rows 72..181 down to 73..182 (matching the rest of the diff, which is
# a pure downward shift of every existing record by one row).
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with the new record's data.
$ws.Cells.Item(72, 1).Value = 7
$ws.Cells.Item(72, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(72, 3).Value = "Ñuble"
$ws.Cells.Item(72, 4).Value = 44579
$ws.Cells.Item(72, 5).Value = 16
$ws.Cells.Item(72, 6).Value = 100112032
$ws.Cells.Item(72, 7).Value = "Zapallo italiano"
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 160
$ws.Cells.Item(72, 11).Value = 9000
$ws.Cells.Item(72, 12).Value = 9500
$ws.Cells.Item(72, 13).Value = 9250
$ws.Cells.Item(72, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(72, 15).Value = "Región del Maule"
$ws.Cells.Item(72, 16).Value = 154
$ws.Cells.Item(72, 17).Value = 60
$ws.Cells.Item(72, 18).Value = "Hortaliza"
